# Apply updated odds values to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value  = 4.1
$ws.Range("I2").Value  = 1.83
$ws.Range("O2").Value  = 1.29
$ws.Range("P2").Value  = 3.75
$ws.Range("Q2").Value  = 1.93
$ws.Range("R2").Value  = 1.93
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 7
$ws.Range("AJ2").Value = 8.5

# Row 4
$ws.Range("G4").Value  = 1.44
$ws.Range("H4").Value  = 4.2
$ws.Range("I4").Value  = 6
$ws.Range("J4").Value  = 2
$ws.Range("K4").Value  = 2.5
$ws.Range("L4").Value  = 6
$ws.Range("Q4").Value  = 1.6
$ws.Range("R4").Value  = 2.3
$ws.Range("AD4").Value = 23
$ws.Range("AG4").Value = 17
$ws.Range("AH4").Value = 51
$ws.Range("AI4").Value = 19
$ws.Range("AK4").Value = 19
$ws.Range("AL4").Value = 67
$ws.Range("AP4").Value = 1.98
$ws.Range("AQ4").Value = 1.83

# Row 7
$ws.Range("K7").Value  = 3.75
$ws.Range("O7").Value  = 1.07
$ws.Range("P7").Value  = 9
$ws.Range("W7").Value  = 2
$ws.Range("X7").Value  = 1.75
$ws.Range("AJ7").Value = 126
$ws.Range("AM7").Value = 151
$ws.Range("AN7").Value = 101

# Row 8
$ws.Range("Q8").Value  = 1.67
$ws.Range("R8").Value  = 2.15
